$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the date text "5.16.17" -> "05.16.17" in column A (rows 2-25).
#    These cells share one sharedString entry, so a single range write updates
#    every row at once. Force text mode first so Excel doesn't reinterpret the
#    dotted string as a serial date.
$ws.Range("A2:A25").NumberFormat = "@"
$ws.Range("A2:A25").Value = "05.16.17"
$ws.Range("A2:A25").NumberFormat = "General"

# 2) Re-apply formatting to the data rows (A3:A25) - this mirrors the
#    original edit's new cell style (same font, explicitly re-applied).
$ws.Range("A3:A25").Style = $ws.Range("A3:A25").Style

# 3) The data rows' height was recalculated (15 -> 13.8) after the edit.
$ws.Range("A2:A25").RowHeight = 13.8

# 4) Leave the selection where the user ended up: A2 active, A2:A25 selected.
[void]$ws.Range("A2:A25").Select()
